$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 used to hold its own unique string ("Row 1 Value" at shared-string index 1,
# distinct from A1's "Row 2 Value"). It now holds the same "Row 1 Value" text as
# A1 -- which collapses the old "Row 2 Value" string out of the shared table.
$ws.Range("A1").Value = "Row 1 Value"
$ws.Range("B1").Value = "Row 1 Value"

# Selection moves from D2 to B1.
$ws.Range("B1").Select() | Out-Null

# Page setup: portrait orientation now recorded for the sheet (xlPortrait = 1).
$ws.PageSetup.Orientation = 1
